$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.916.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.09%  "
$ws.Range("D3").Value = "'3.074.65"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.41%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'579.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.64%  "
$ws.Range("D6").Value = "'141.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.15%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'3.064.84"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.50%  "
$ws.Range("E9").Value = "  +0.82%  "
$ws.Range("E10").Value = "  +4.73%  "
$ws.Range("D11").Value = "'5.76"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +11.53%  "
$ws.Range("E12").Value = "  +1.92%  "
$ws.Range("D13").Value = "'0.0000239"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.69%  "
$ws.Range("D14").Value = "'35.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.12%  "
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("D16").Value = "'3.584.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.48%  "
$ws.Range("E17").Value = "  -0.94%  "
$ws.Range("D18").Value = "'3.073.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.59%  "
$ws.Range("D19").Value = "'61.847.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.23%  "
$ws.Range("D20").Value = "'447.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.64%  "
$ws.Range("D21").Value = "'13.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.70%  "
$ws.Range("E22").Value = "  +1.27%  "
$ws.Range("D23").Value = "'7.43"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.84%  "
$ws.Range("D24").Value = "'13.75"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.46%  "
$ws.Range("D25").Value = "'81.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.89%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  +3.95%  "
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'2.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.31%  "
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "'8.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.27%  "
$ws.Range("D31").Value = "'6.77"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.12%  "
$ws.Range("E32").Value = "  +12.26%  "
$ws.Range("D33").Value = "'26.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.06%  "
$ws.Range("E34").Value = "  +4.70%  "
$ws.Range("D35").Value = "'0.0₃0792"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.75%  "
$ws.Range("E36").Value = "  +2.19%  "
$ws.Range("E37").Value = "  +4.08%  "
$ws.Range("D38").Value = "'50.06"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.82%  "
$ws.Range("D39").Value = "'2.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.73%  "
$ws.Range("D40").Value = "'8.78"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.84%  "
$ws.Range("D41").Value = "'420.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.39%  "
$ws.Range("D42").Value = "'2.944.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.23%  "
$ws.Range("E43").Value = "  +5.13%  "
$ws.Range("D44").Value = "'0.275"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.99%  "
$ws.Range("E45").Value = "  +0.40%  "
$ws.Range("D46").Value = "'2.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.85%  "
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "'0.999"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "'124.35"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.11%  "
$ws.Range("B49").Value = "Arweave"
$ws.Range("C49").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D49").Value = "'35.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("D51").Value = "'24.33"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.14%  "
